$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 515.6667
$ws.Range("I12").Value = 448.85715
$ws.Range("K12").Value = 448.85715
$ws.Range("M12").Value = -278.85715
$ws.Range("H15").Value = 2434.0417
$ws.Range("I15").Value = 2434.0417
$ws.Range("K15").Value = 7302.125100000001
$ws.Range("M15").Value = -7133.125100000001
$ws.Range("H32").Value = 8566.333000000001
$ws.Range("J32").Value = 8566.333000000001
$ws.Range("L32").Value = 8566.333000000001
$ws.Range("N32").Value = -9218.333000000001
$ws.Range("H113").Value = 15931.833
$ws.Range("I113").Value = 18150.334
$ws.Range("K113").Value = 18150.334
$ws.Range("M113").Value = -14896.334
$ws.Range("H115").Value = 566.7143
$ws.Range("I115").Value = 454.375
$ws.Range("K115").Value = 1363.125
$ws.Range("M115").Value = 203.875
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = $null
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1100.1837
$ws.Range("I2").Value = 841.89655
$ws.Range("K2").Value = 841.89655
$ws.Range("M2").Value = -728.89655
$ws.Range("H10").Value = 10002
$ws.Range("I10").Value = 10002
$ws.Range("K10").Value = 10002
$ws.Range("M10").Value = -9832
$ws.Range("H32").Value = 9854.593000000001
$ws.Range("I32").Value = 8207.291999999999
$ws.Range("K32").Value = 8207.291999999999
$ws.Range("M32").Value = -7920.291999999999
$ws.Range("H45").Value = 2546.111
$ws.Range("I45").Value = 2259.6365
$ws.Range("K45").Value = 2259.6365
$ws.Range("M45").Value = -1882.6365
$ws.Range("H74").Value = 14174.576
$ws.Range("I74").Value = 1291
$ws.Range("J74").Value = 21536.62
$ws.Range("K74").Value = 1291
$ws.Range("L74").Value = 21536.62
$ws.Range("M74").Value = -417
$ws.Range("N74").Value = -23284.62
$ws.Range("H77").Value = 14174.576
$ws.Range("I77").Value = 1291
$ws.Range("J77").Value = 21536.62
$ws.Range("K77").Value = 6455
$ws.Range("L77").Value = 107683.1
$ws.Range("M77").Value = -2087
$ws.Range("N77").Value = -116419.1
$ws.Range("H110").Value = 1952.579
$ws.Range("I110").Value = 2266.4138
$ws.Range("J110").Value = 941.3333
$ws.Range("K110").Value = 2266.4138
$ws.Range("L110").Value = 941.3333
$ws.Range("M110").Value = -221.4137999999998
$ws.Range("N110").Value = -5031.3333
$ws.Range("H116").Value = 1100.1837
$ws.Range("I116").Value = 841.89655
$ws.Range("K116").Value = 841.89655
$ws.Range("M116").Value = 1452.10345
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1100.1837
$ws.Range("I3").Value = 841.89655
$ws.Range("K3").Value = 841.89655
$ws.Range("M3").Value = -727.89655
$ws.Range("H20").Value = 17996.334
$ws.Range("I20").Value = 5750.64
$ws.Range("K20").Value = 5750.64
$ws.Range("M20").Value = -5503.64
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 250118.25
$ws.Range("I10").Value = 197
$ws.Range("J10").Value = 500039.5
$ws.Range("K10").Value = 197
$ws.Range("L10").Value = 500039.5
$ws.Range("M10").Value = -58
$ws.Range("N10").Value = -500317.5
$ws.Range("H16").Value = 2286.25
$ws.Range("J16").Value = 2283.3333
$ws.Range("L16").Value = 2283.3333
$ws.Range("N16").Value = -2857.3333
$ws.Range("H25").Value = 1200
$ws.Range("I25").Value = 1800
$ws.Range("K25").Value = 1800
$ws.Range("M25").Value = -1626
$ws.Range("H39").Value = 6525
$ws.Range("I39").Value = 6525
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6525
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -6134
$ws.Range("N39").Value = $null
$ws.Range("H49").Value = 6525
$ws.Range("I49").Value = 6525
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 6525
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -6343
$ws.Range("N49").Value = $null
$ws.Range("H58").Value = 15401.2
$ws.Range("J58").Value = 19984.79
$ws.Range("L58").Value = 19984.79
$ws.Range("N58").Value = -20390.79
$ws.Range("H107").Value = 2300
$ws.Range("I107").Value = 3541.4
$ws.Range("J107").Value = 1265.5
$ws.Range("K107").Value = 3541.4
$ws.Range("L107").Value = 1265.5
$ws.Range("M107").Value = -1621.4
$ws.Range("N107").Value = -5105.5
$ws.Range("H113").Value = 2286.25
$ws.Range("J113").Value = 2283.3333
$ws.Range("L113").Value = 2283.3333
$ws.Range("N113").Value = -6623.3333
$ws.Range("H136").Value = 15401.2
$ws.Range("J136").Value = 19984.79
$ws.Range("M136").Value = -30071.858
$ws.Range("N136").Value = -65054.37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2711147
$ws.Range("I5").Value = 1260
$ws.Range("J5").Value = 24390244
$ws.Range("K5").Value = 3780
$ws.Range("L5").Value = 73170732
$ws.Range("M5").Value = -3668
$ws.Range("N5").Value = -73170956
$ws.Range("H14").Value = 207.4
$ws.Range("I14").Value = 207.4
$ws.Range("K14").Value = 622.2
$ws.Range("M14").Value = -449.2
$ws.Range("H96").Value = 2092.8572
$ws.Range("I96").Value = 2025
$ws.Range("K96").Value = 6075
$ws.Range("M96").Value = -4016
$ws.Range("H113").Value = 1793.7142
$ws.Range("I113").Value = 1182.6666
$ws.Range("J113").Value = 2252
$ws.Range("K113").Value = 3547.9998
$ws.Range("L113").Value = 6756
$ws.Range("M113").Value = -1377.9998
$ws.Range("N113").Value = -11096
$ws.Range("H131").Value = 1488
$ws.Range("I131").Value = 1340.1666
$ws.Range("J131").Value = 1497.4362
$ws.Range("K131").Value = 4020.4998
$ws.Range("L131").Value = 4492.3086
$ws.Range("M131").Value = 1019.5002
$ws.Range("N131").Value = -14572.3086
$ws.Range("H135").Value = 2711147
$ws.Range("I135").Value = 1260
$ws.Range("J135").Value = 24390244
$ws.Range("K135").Value = 11340
$ws.Range("L135").Value = 219512196
$ws.Range("M135").Value = -8805
$ws.Range("N135").Value = -219517266
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8500.388999999999
$ws.Range("J126").Value = 6840.727
$ws.Range("L126").Value = 20522.181
$ws.Range("N126").Value = -25462.181
$ws.Range("H132").Value = 15244.294
$ws.Range("I132").Value = 10867.286
$ws.Range("K132").Value = 32601.858
$ws.Range("M132").Value = -30071.858
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10226
$ws.Range("H22").Value = 4610.08
$ws.Range("I22").Value = 2244.0527
$ws.Range("J22").Value = 12102.5
$ws.Range("K22").Value = 2244.0527
$ws.Range("L22").Value = 12102.5
$ws.Range("M22").Value = -1949.0527
$ws.Range("N22").Value = -12692.5
$ws.Range("H25").Value = 158108.08
$ws.Range("I25").Value = 1905.0952
$ws.Range("J25").Value = 978173.75
$ws.Range("K25").Value = 1905.0952
$ws.Range("L25").Value = 978173.75
$ws.Range("M25").Value = -1675.0952
$ws.Range("N25").Value = -978633.75
$ws.Range("H27").Value = 4610.08
$ws.Range("I27").Value = 2244.0527
$ws.Range("J27").Value = 12102.5
$ws.Range("K27").Value = 2244.0527
$ws.Range("L27").Value = 12102.5
$ws.Range("M27").Value = -2137.0527
$ws.Range("N27").Value = -12316.5
$ws.Range("H28").Value = 10000
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10464
$ws.Range("H37").Value = 10000
$ws.Range("J37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("N37").Value = -10214
$ws.Range("H40").Value = 5529.5
$ws.Range("J40").Value = 5650.9165
$ws.Range("L40").Value = 5650.9165
$ws.Range("N40").Value = -5922.9165
$ws.Range("H46").Value = 3177.6316
$ws.Range("I46").Value = 1781.1666
$ws.Range("K46").Value = 1781.1666
$ws.Range("M46").Value = -1593.1666
$ws.Range("H132").Value = 1492007.5
$ws.Range("I132").Value = 1982.5
$ws.Range("J132").Value = 2684027.5
$ws.Range("K132").Value = 5947.5
$ws.Range("L132").Value = 8052082.5
$ws.Range("M132").Value = -3417.5
$ws.Range("N132").Value = -8057142.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4500
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5344
